# Update the "Förändrad" (changed/updated) date column (C) from 45183 to 45184
# for all data rows (2 through 12) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45183) {
        $cell.Value = 45184
    }
}
